$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AC (28 columns after A) hold the match data for each row.
# Swap the full row content (except column A, the sequential id) between
# the given row pairs.
function Swap-Rows($ws, $r1, $r2) {
    $rng1 = $ws.Range("B$r1`:AC$r1")
    $rng2 = $ws.Range("B$r2`:AC$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

Swap-Rows $ws 123 124
Swap-Rows $ws 131 132
Swap-Rows $ws 144 145
Swap-Rows $ws 191 192
Swap-Rows $ws 202 203

# Individual odds updates for rows 210-214
$ws.Range("N210").Value = 5
$ws.Range("P210").Value = 1.6
$ws.Range("U210").Value = 1.925
$ws.Range("V210").Value = 1.925

$ws.Range("R211").Value = 1.775
$ws.Range("S211").Value = 2.1

$ws.Range("N212").Value = 2.375
$ws.Range("O212").Value = 3.1
$ws.Range("P212").Value = 3.1
$ws.Range("R212").Value = 2.05
$ws.Range("S212").Value = 1.8
$ws.Range("T212").Value = 2.5

$ws.Range("R213").Value = 1.875
$ws.Range("S213").Value = 1.975

$ws.Range("N214").Value = 2.15
$ws.Range("O214").Value = 3.25
$ws.Range("T214").Value = 2.5
$ws.Range("U214").Value = 2
$ws.Range("V214").Value = 1.85
